$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 11851
$ws.Range("C3").Value = 6244
$ws.Range("C4").Value = 9174
$ws.Range("C5").Value = 6907
$ws.Range("C6").Value = 5297
$ws.Range("C7").Value = 8725
$ws.Range("C8").Value = 23062
$ws.Range("C9").Value = 17104
$ws.Range("C10").Value = 5155
$ws.Range("C11").Value = 3959
$ws.Range("C12").Value = 54
